# New daily allocation row appended after running the profit-split job on 2025-09-27.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 26

# Force the date-looking label to stay a text string (matches the existing
# "Date" column cells, which are plain text, not real dates) instead of
# letting Excel auto-convert "09/27/2025" into a date serial on entry.
$ws.Cells.Item($row, 1).NumberFormat = "@"
$ws.Cells.Item($row, 1).Value = "09/27/2025"
# Drop back to the default "Normal" cell style so the new row matches the
# unstyled look of the other data rows (no leftover text-format style).
$ws.Cells.Item($row, 1).Style = "Normal"

$ws.Cells.Item($row, 2).Value = 0.1324659859615432
$ws.Cells.Item($row, 3).Value = 0.8675340140384568
